# Updated symbol list on Fri Feb 17 15:23:34 UTC 2023 with GitHub Actions
#
# Refreshes the "cryptos" price/volume snapshot on Sheet1: updates Price
# (column D) and Volume(1h) (column E) for the rows whose underlying coin
# data moved since the last run, and re-syncs the "Coin"/"Link" columns
# (B/C) for rows 15-23 where the coin ranking list shifted by one
# position (the coin previously in the next row moves up, and the coin
# that fell out of the top of that block reappears at the bottom).
#
# Columns D and E hold plain text (e.g. "309.75", "-2.95%"), not numbers,
# so each numeric-looking value is written with a leading apostrophe to
# force Excel to keep storing it as text instead of auto-converting it to
# a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.75"
$ws.Range("E2").Value = "'-2.95%"
$ws.Range("D3").Value = "'53.87"
$ws.Range("E3").Value = "'8.45%"
$ws.Range("D4").Value = "'5.108"
$ws.Range("E4").Value = "'-2.94%"
$ws.Range("D5").Value = "'0.07847"
$ws.Range("E5").Value = "'-1.23%"
$ws.Range("E6").Value = "'-1.35%"
$ws.Range("D7").Value = "'1.366"
$ws.Range("E7").Value = "'-0.96%"
$ws.Range("D8").Value = "'1.567"
$ws.Range("E8").Value = "'-4.30%"
$ws.Range("D9").Value = "'0.1220"
$ws.Range("E9").Value = "'-5.37%"
$ws.Range("D10").Value = "'0.2003"
$ws.Range("E10").Value = "'1.95%"
$ws.Range("D11").Value = "'0.04722"
$ws.Range("E11").Value = "'2.65%"
$ws.Range("D12").Value = "'0.09467"
$ws.Range("E12").Value = "'0.46%"
$ws.Range("E13").Value = "'-0.27%"
$ws.Range("D14").Value = "'0.001261"
$ws.Range("E14").Value = "'-4.22%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005798"
$ws.Range("E15").Value = "'-2.17%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007487"
$ws.Range("E16").Value = "'2,017.75%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.331"
$ws.Range("E17").Value = "'-0.38%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.419"
$ws.Range("E18").Value = "'-0.61%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3445"
$ws.Range("E19").Value = "'-0.34%"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'8.025"
$ws.Range("E20").Value = "'-2.16%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1366"
$ws.Range("E21").Value = "'-1.67%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.3082"
$ws.Range("E22").Value = "'-0.28%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04167"
$ws.Range("E23").Value = "'-0.03%"
$ws.Range("D24").Value = "'0.001256"
$ws.Range("E24").Value = "'-4.47%"
$ws.Range("D25").Value = "'0.003916"
$ws.Range("E25").Value = "'-7.99%"
$ws.Range("D26").Value = "'0.0001346"
$ws.Range("E26").Value = "'-0.17%"
$ws.Range("D38").Value = "'0.02608"
$ws.Range("E38").Value = "'-3.09%"
$ws.Range("D39").Value = "'0.05943"
$ws.Range("E39").Value = "'3.09%"
$ws.Range("D40").Value = "'0.01052"
$ws.Range("E40").Value = "'-3.57%"
$ws.Range("D41").Value = "'0.007881"
$ws.Range("E41").Value = "'-1.56%"
$ws.Range("D42").Value = "'0.1421"
$ws.Range("E42").Value = "'-1.31%"
$ws.Range("D43").Value = "'0.008183"
$ws.Range("E43").Value = "'6.49%"
$ws.Range("D44").Value = "'0.008453"
$ws.Range("E44").Value = "'-0.11%"
$ws.Range("D45").Value = "'0.3114"
$ws.Range("E45").Value = "'-2.38%"
$ws.Range("D46").Value = "'0.00007243"
$ws.Range("E46").Value = "'9.28%"
$ws.Range("D47").Value = "'0.00000000747"
$ws.Range("E47").Value = "'-0.27%"
$ws.Range("D48").Value = "'0.05578"
$ws.Range("E48").Value = "'1.43%"
$ws.Range("D49").Value = "'0.002609"
$ws.Range("E49").Value = "'-34.71%"
$ws.Range("D50").Value = "'0.00002091"
$ws.Range("E50").Value = "'-0.27%"
$ws.Range("D51").Value = "'0.0001992"
$ws.Range("E51").Value = "'-0.27%"
